$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = -12.74230000000001
$ws.Range("C12").Value = -14.38810000000001
$ws.Range("E13").Value = 12.0817
$ws.Range("C18").Value = -14.2532
